# The sheet holds one weekly price record per row (rows 2..143, header on
# row 1). This commit adds one new weekly record. In the canonical XML the
# new record lands at row 33, which pushes every existing row from the old
# row 33 down by one (old row 33 -> new row 34, ..., old row 143 -> new row
# 144), growing the used range from A1:R143 to A1:R144.
#
# Reproduce that with a real row insert (shifting existing rows down, and
# carrying their formatting with them) and then populate the freshly
# inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 33; rows 33:143 shift down to 34:144.
$ws.Rows("33:33").Insert()

# Fill the new row 33 with the new weekly Ajo price record.
$ws.Range("A33").Value = 8
$ws.Range("B33").Value = "Terminal La Palmera de La Serena"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = 44459
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 100112003
$ws.Range("G33").Value = "Ajo"
$ws.Range("H33").Value = "Chino"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 16000
$ws.Range("M33").Value = 15500
$ws.Range("N33").Value = "`$/caja 10 kilos"
$ws.Range("O33").Value = "China"
$ws.Range("P33").Value = 1550
$ws.Range("Q33").Value = 10
$ws.Range("R33").Value = "Hortaliza"
